$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.178.42"
$ws.Range("E2").Value = "  +2.09%  "
$ws.Range("D3").Value = "3.629.37"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'195.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.17%  "
$ws.Range("D6").Value = "'577.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("D7").Value = "3.624.24"
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("E8").Value = "  +3.36%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("D10").Value = "'0.677"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.11%  "
$ws.Range("D11").Value = "'0.153"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +8.36%  "
$ws.Range("D12").Value = "'55.87"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.00%  "
$ws.Range("D13").Value = "'0.0000295"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +22.55%  "
$ws.Range("D14").Value = "'10.12"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.54%  "
$ws.Range("D15").Value = "4.207.43"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "3.634.47"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").Value = "'12.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +4.17%  "
$ws.Range("D19").Value = "68.121.37"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Value = "'18.53"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.29%  "
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").Value = "'403.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.67%  "
$ws.Range("D23").Value = "'12.76"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +25.20%  "
$ws.Range("D24").Value = "'4.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'86.03"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("E26").Value = "  +5.44%  "
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("D28").Value = "'3.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.10%  "
$ws.Range("D29").Value = "'6.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "'8.12"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +22.28%  "
$ws.Range("D31").Value = "'9.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.30%  "
$ws.Range("D32").Value = "'31.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("D33").Value = "'688.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +19.10%  "
$ws.Range("D34").Value = "'12.21"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.56%  "
$ws.Range("E35").Value = "  +7.86%  "
$ws.Range("D36").Value = "'64.81"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").Value = "'42.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.49%  "
$ws.Range("D38").Value = "'0.416"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.76%  "
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "0.0₃0797"
$ws.Range("E40").Value = "  +11.97%  "
$ws.Range("E41").Value = "  +25.02%  "
$ws.Range("E42").Value = "  +15.22%  "
$ws.Range("E43").Value = "  +3.53%  "
$ws.Range("D44").Value = "3.159.11"
$ws.Range("E44").Value = "  +19.02%  "
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").Value = "'2.86"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +29.21%  "
$ws.Range("D47").Value = "'0.0424"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("D48").Value = "'0.132"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.83%  "
$ws.Range("D49").Value = "'8.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.38%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").Value = "'142.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.22%  "
$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "'3.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.40%  "
